# Applies:
#  - Bugfix: job "Process.time" column values (B) corrected to 25
#  - Bugfix: "Due Date (Seconds)" column values (G) corrected to 22
#  - New columns: "Penalty Rate" (H, value 1) and "Cost" (I, value 2)
#  - Table1 grows to include the two new columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing data: Process.time (B) and Due Date (Seconds) (G) ---
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 25
    $ws.Cells.Item($r, 7).Value = 22
}

# --- Grow the table to include the new columns ---
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:I11"))

# --- Add the two new header cells (must happen after Resize, which stamps
#     default "ColumnN" headers into the newly added range) ---
$ws.Cells.Item(1, 8).Value = "Penalty Rate"
$ws.Cells.Item(1, 9).Value = "Cost"

# --- Populate the new columns for every data row ---
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
    $ws.Cells.Item($r, 9).Value = 2
}

# --- Update selection to match the final state ---
$ws.Range("I11").Select()
